# DATA: update 2020-04-03
#
# Append the newest daily COVID-19 Indonesia figures as a new row at the
# bottom of the Sheet1 list (row 46), then move the view/selection down to
# follow the newly entered data - mirroring what a person does in Excel
# when they type a new row under an existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 46
$lastRow = $newRow - 1

# Copy the formatting of the previous last row down into the new row so
# the date cell (column A) keeps its existing "yyyy-mm-dd" number format
# and the rest of the row matches the table's normal (unstriped) style.
$ws.Range("A" + $lastRow + ":H" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":H" + $newRow).PasteSpecial()
$excel.CutCopyMode = $false

# tanggal, sembuh, meninggal, negatif, kasus_perawatan, proses_periksa,
# jumlah_periksa, konfirmasi
$ws.Cells.Item($newRow, 1).Value = 43923
$ws.Cells.Item($newRow, 2).Value = 7425
$ws.Cells.Item($newRow, 3).Value = 1790
$ws.Cells.Item($newRow, 4).Value = 112
$ws.Cells.Item($newRow, 5).Value = 170
$ws.Cells.Item($newRow, 6).Value = 5635
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0

# Scroll the window down a few rows and select the new bottom-right cell,
# matching the saved view (topLeftCell A22, selection I46).
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("I" + $newRow).Select()
